# UndoRedoStartingStackDiagram.pptx
# "DeveloperGuide: Revise some of the figures to match color scheme"
#
# The only content edit in this single-slide deck that is reachable
# through the PowerPoint object model is the renaming of the sample
# variable shown in the "Table 20" diagram cell from `prevAddressBook`
# to `prevCoinBook` (the rest of the stored `undo`/`redo` example value
# is left untouched).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the table shape ("Table 20") on the slide instead of hard-coding
# its index, in case shape ordering ever shifts.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

$tbl = $tableShape.Table

# The second row holds the two-line cell:
#   targetIndex = 5
#   prevAddressBook = s3
$cell = $tbl.Cell(2, 1)
$tr = $cell.Shape.TextFrame.TextRange
$paragraphs = $tr.Paragraphs()

for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    if ($para.Text -match "prevAddressBook") {
        $para.Text = $para.Text -replace "prevAddressBook", "prevCoinBook"
    }
}
